$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. "583.68") are not auto-converted to numbers by Excel, then
# restore the default "Normal" style so no stray style index is left
# on the cells (matches original workbook which has no s= attribute).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '67.773.40'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '3.327.74'
$ws.Range("E3").Value = '  +1.18%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '583.68'
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").Value = '174.52'
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '0.589'
$ws.Range("E8").Value = '  +1.72%  '
$ws.Range("D9").Value = '3.323.25'
$ws.Range("E9").Value = '  +1.18%  '
$ws.Range("D10").Value = '0.181'
$ws.Range("E10").Value = '  +4.47%  '
$ws.Range("D11").Value = '0.578'
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("D12").Value = '46.97'
$ws.Range("E12").Value = '  +3.99%  '
$ws.Range("D13").Value = '0.0000272'
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").Value = '698.93'
$ws.Range("E14").Value = '  +5.06%  '
$ws.Range("D15").Value = '3.865.13'
$ws.Range("E15").Value = '  +1.15%  '
$ws.Range("D16").Value = '8.36'
$ws.Range("E16").Value = '  +0.71%  '
$ws.Range("D17").Value = '67.759.82'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '0.119'
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").Value = '3.326.54'
$ws.Range("E19").Value = '  +0.99%  '
$ws.Range("D20").Value = '17.45'
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("D21").Value = '11.11'
$ws.Range("E21").Value = '  +2.60%  '
$ws.Range("D22").Value = '0.888'
$ws.Range("E22").Value = '  +0.84%  '
$ws.Range("D23").Value = '5.41'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").Value = '16.91'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").Value = '101.29'
$ws.Range("E25").Value = '  +3.11%  '
$ws.Range("D26").Value = '3.90'
$ws.Range("E26").Value = '  +1.56%  '
$ws.Range("D27").Value = '2.68'
$ws.Range("E27").Value = '  +1.36%  '
$ws.Range("D28").Value = '9.41'
$ws.Range("E28").Value = '  +3.01%  '
$ws.Range("D29").Value = '32.83'
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("D30").Value = '8.52'
$ws.Range("E30").Value = '  +2.27%  '
$ws.Range("D31").Value = '6.96'
$ws.Range("E31").Value = '  -0.36%  '
$ws.Range("D32").Value = '573.63'
$ws.Range("E32").Value = '  -0.88%  '
$ws.Range("D33").Value = '10.99'
$ws.Range("E33").Value = '  +0.93%  '
$ws.Range("D34").Value = '0.105'
$ws.Range("E34").Value = '  +2.32%  '
$ws.Range("D35").Value = '3.725.68'
$ws.Range("E35").Value = '  -0.45%  '
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("D37").Value = '56.48'
$ws.Range("E37").Value = '  +1.63%  '
$ws.Range("D38").Value = '3.27'
$ws.Range("E38").Value = '  -2.49%  '
$ws.Range("D39").Value = '35.54'
$ws.Range("E39").Value = '  +10.27%  '
$ws.Range("D40").Value = '0.134'
$ws.Range("E40").Value = '  +2.65%  '
$ws.Range("D41").Value = '3.13'
$ws.Range("E41").Value = '  +3.16%  '
$ws.Range("D42").Value = '2.60'
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("D43").Value = '0.0₃0671'
$ws.Range("E43").Value = '  +1.48%  '
$ws.Range("D44").Value = '0.333'
$ws.Range("E44").Value = '  +2.35%  '
$ws.Range("D45").Value = '3.30'
$ws.Range("E45").Value = '  +1.79%  '
$ws.Range("D46").Value = '0.0405'
$ws.Range("E46").Value = '  +1.21%  '
$ws.Range("D47").Value = '2.61'
$ws.Range("E47").Value = '  +1.58%  '
$ws.Range("D48").Value = '0.128'
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").Value = '1.32'
$ws.Range("E50").Value = '  -2.26%  '
$ws.Range("D51").Value = '130.78'
$ws.Range("E51").Value = '  +1.58%  '

$dRange.Style = "Normal"

